$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B8").Value = "SingleUseId7"
$ws.Range("C8").Value = "OurTypography"
$ws.Range("D8").Value = "Left"
$ws.Range("E8").Value = "essa"
$ws.Range("F8").Value = "LTR"
